$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Demon', ['Token Creature — Demon', 'Flying', '*/*'])"
$ws.Range("A3").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
$ws.Range("A4").Value = "('Thrull', ['Token Creature — Thrull', '0/1'])"

$ws.Range("A5:A12").Clear()
